$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AI5").Value = 17
$ws.Range("AN5").Value = 8.5
$ws.Range("G5").Value = 5.5
$ws.Range("H5").Value = 4.33
$ws.Range("I5").Value = 1.57
$ws.Range("Y5").Value = 1.8
$ws.Range("Z5").Value = 1.95

$ws.Range("AB12").Value = 6.5
$ws.Range("AF12").Value = 29
$ws.Range("AJ12").Value = 67
$ws.Range("G12").Value = 1.45
$ws.Range("I12").Value = 8
$ws.Range("L12").Value = 7
$ws.Range("N12").Value = 10
$ws.Range("Q12").Value = 2.02
$ws.Range("R12").Value = 1.88
$ws.Range("U12").Value = 3.4
$ws.Range("V12").Value = 1.33

$ws.Range("AL21").Value = 10
$ws.Range("AN21").Value = 8.5
$ws.Range("G21").Value = 6
$ws.Range("I21").Value = 1.45
$ws.Range("N21").Value = 19
$ws.Range("O21").Value = 1.14
$ws.Range("P21").Value = 5.5
$ws.Range("W21").Value = 1.25
$ws.Range("X21").Value = 3.75

$ws.Range("AB22").Value = 11
$ws.Range("AC22").Value = 10
$ws.Range("AD22").Value = 23
$ws.Range("AE22").Value = 23
$ws.Range("AL22").Value = 7.5
$ws.Range("AM22").Value = 13
$ws.Range("AR22").Value = 1.8
$ws.Range("AS22").Value = 2.05
$ws.Range("G22").Value = 2.55
$ws.Range("I22").Value = 2.9
$ws.Range("J22").Value = 3.25
$ws.Range("L22").Value = 3.75
$ws.Range("M22").Value = 1.08
$ws.Range("N22").Value = 8
$ws.Range("Q22").Value = 2.35
$ws.Range("R22").Value = 1.57
$ws.Range("S22").Value = 3.85

$ws.Range("AF23").Value = 67
$ws.Range("H23").Value = 2.7
$ws.Range("O23").Value = 1.83
$ws.Range("P23").Value = 1.83
$ws.Range("S23").Value = 8
$ws.Range("W23").Value = 1.9
$ws.Range("X23").Value = 1.9
$ws.Range("Y23").Value = 3
$ws.Range("Z23").Value = 1.36

$ws.Range("AD24").Value = 15
$ws.Range("AG24").Value = 7
$ws.Range("AL24").Value = 9
$ws.Range("G24").Value = 1.9
$ws.Range("H24").Value = 3.3
$ws.Range("I24").Value = 4.33
$ws.Range("J24").Value = 2.63
$ws.Range("K24").Value = 1.95
$ws.Range("L24").Value = 5.5
$ws.Range("Y24").Value = 2.2
$ws.Range("Z24").Value = 1.62

$ws.Range("M26").Value = 1.14
$ws.Range("N26").Value = 5.5

$ws.Range("AB27").Value = 7.5
$ws.Range("AE27").Value = 26
$ws.Range("AG27").Value = 4.5
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 2.88
$ws.Range("I27").Value = 4.5
$ws.Range("J27").Value = 2.88
$ws.Range("O27").Value = 1.73
$ws.Range("P27").Value = 2
$ws.Range("W27").Value = 1.78
$ws.Range("X27").Value = 2.03

$ws.Range("AB28").Value = 7.5
$ws.Range("AD28").Value = 17
$ws.Range("AH28").Value = 6.5
$ws.Range("AL28").Value = 8
$ws.Range("AM28").Value = 19
$ws.Range("AN28").Value = 15
$ws.Range("AR28").Value = 2.1
$ws.Range("AS28").Value = 1.78
$ws.Range("G28").Value = 1.95
$ws.Range("H28").Value = 3.2
$ws.Range("I28").Value = 4.1
$ws.Range("J28").Value = 2.75
$ws.Range("L28").Value = 5.5
$ws.Range("M28").Value = 1.11
$ws.Range("N28").Value = 6.5
$ws.Range("O28").Value = 1.57
$ws.Range("P28").Value = 2.25
$ws.Range("Q28").Value = 2.7
$ws.Range("R28").Value = 1.44
$ws.Range("U28").Value = 6
$ws.Range("V28").Value = 1.13
$ws.Range("Y28").Value = 2.38
$ws.Range("Z28").Value = 1.53

$ws.Range("AB29").Value = 9
$ws.Range("AD29").Value = 23
$ws.Range("AO29").Value = 41
$ws.Range("G29").Value = 2.3
$ws.Range("I29").Value = 3.9
$ws.Range("M29").Value = 1.18
$ws.Range("N29").Value = 4.5

$ws.Range("AB31").Value = 8.5
$ws.Range("AC31").Value = 11
$ws.Range("AG31").Value = 5.5
$ws.Range("G31").Value = 2.15
$ws.Range("H31").Value = 3.1
$ws.Range("I31").Value = 3.6
$ws.Range("J31").Value = 3.1
$ws.Range("L31").Value = 4.75
$ws.Range("M31").Value = 1.14
$ws.Range("N31").Value = 5.5

$ws.Range("AB49").Value = 8
$ws.Range("AC49").Value = 9.5
$ws.Range("AD49").Value = 17
$ws.Range("AG49").Value = 6
$ws.Range("AI49").Value = 21
$ws.Range("AM49").Value = 19
$ws.Range("G49").Value = 2
$ws.Range("I49").Value = 4.1
$ws.Range("L49").Value = 4.75
$ws.Range("S49").Value = 4.2
$ws.Range("T49").Value = 1.21

$ws.Range("G50").Value = 2.2
$ws.Range("K50").Value = 1.92
$ws.Range("M50").Value = 1.1
$ws.Range("N50").Value = 7
$ws.Range("R50").Value = 1.47

$ws.Range("AR51").Value = 1.54
$ws.Range("AS51").Value = 2.44
$ws.Range("G51").Value = 1.42
$ws.Range("R51").Value = 1.72
$ws.Range("S51").Value = 2.9
$ws.Range("T51").Value = 1.4

$ws.Range("AA84").Value = 6
$ws.Range("AB84").Value = 6.5
$ws.Range("AG84").Value = 9.5
$ws.Range("AI84").Value = 21
$ws.Range("AK84").Value = 501
$ws.Range("AN84").Value = 19
$ws.Range("AP84").Value = 51
$ws.Range("G84").Value = 1.55
$ws.Range("I84").Value = 6.25
$ws.Range("J84").Value = 2.1
$ws.Range("L84").Value = 6.5
$ws.Range("M84").Value = 1.06
$ws.Range("N84").Value = 10
$ws.Range("O84").Value = 1.33
$ws.Range("P84").Value = 3.25
$ws.Range("Q84").Value = 2.05
$ws.Range("R84").Value = 1.75
$ws.Range("U84").Value = 3.5
$ws.Range("V84").Value = 1.29
$ws.Range("W84").Value = 1.44
$ws.Range("X84").Value = 2.63
$ws.Range("Y84").Value = 2.05
$ws.Range("Z84").Value = 1.7

$ws.Range("AA93").Value = 13
$ws.Range("AG93").Value = 6.5
$ws.Range("AH93").Value = 7
$ws.Range("AJ93").Value = 81
$ws.Range("AP93").Value = 17
$ws.Range("AR93").Value = 1.8
$ws.Range("AS93").Value = 2.05
$ws.Range("G93").Value = 7.5
$ws.Range("H93").Value = 3.4
$ws.Range("I93").Value = 1.57
$ws.Range("J93").Value = 7
$ws.Range("K93").Value = 2.05
$ws.Range("L93").Value = 2.2
$ws.Range("M93").Value = 1.1
$ws.Range("N93").Value = 7
$ws.Range("O93").Value = 1.44
$ws.Range("P93").Value = 2.63
$ws.Range("Q93").Value = 2.4
$ws.Range("R93").Value = 1.53
$ws.Range("S93").Value = 3.7
$ws.Range("T93").Value = 1.28
$ws.Range("U93").Value = 4.5
$ws.Range("V93").Value = 1.18
$ws.Range("W93").Value = 1.53
$ws.Range("X93").Value = 2.38

$ws.Range("AR94").Value = 1.56
$ws.Range("AS94").Value = 2.46
$ws.Range("M94").Value = 1.06
$ws.Range("N94").Value = 10
$ws.Range("S94").Value = 2.95
$ws.Range("T94").Value = 1.41
$ws.Range("U94").Value = 3.5
$ws.Range("V94").Value = 1.29

$ws.Range("AE101").Value = 19
$ws.Range("AG101").Value = 7
$ws.Range("G101").Value = 1.9
$ws.Range("H101").Value = 3.25
$ws.Range("J101").Value = 2.63

$ws.Range("AA137").Value = 5.6
$ws.Range("AB137").Value = 7
$ws.Range("AC137").Value = 8.25
$ws.Range("AG137").Value = 7.6
$ws.Range("AH137").Value = 6.6
$ws.Range("AL137").Value = 11
$ws.Range("AM137").Value = 28
$ws.Range("AO137").Value = 110
$ws.Range("AQ137").Value = 75
$ws.Range("G137").Value = 1.7
$ws.Range("I137").Value = 5
$ws.Range("K137").Value = 2.05
$ws.Range("L137").Value = 5.4
$ws.Range("P137").Value = 2.57
$ws.Range("W137").Value = 1.44
$ws.Range("X137").Value = 2.42

